$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81; this shifts existing rows 81..153 down to 82..154,
# matching the dimension growing from A1:R153 to A1:R154.
$ws.Rows("81").Insert()

# Populate the newly inserted row 81 with the new weekly record.
$ws.Cells.Item(81, 1).Value = 11
$ws.Cells.Item(81, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(81, 3).Value = "Bíobío"
$ws.Cells.Item(81, 4).Value = 45049
$ws.Cells.Item(81, 5).Value = 8
$ws.Cells.Item(81, 6).Value = 100112001
$ws.Cells.Item(81, 7).Value = "Berenjena"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 100
$ws.Cells.Item(81, 11).Value = 7500
$ws.Cells.Item(81, 12).Value = 8500
$ws.Cells.Item(81, 13).Value = 8000
$ws.Cells.Item(81, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(81, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(81, 16).Value = 133
$ws.Cells.Item(81, 17).Value = 60
$ws.Cells.Item(81, 18).Value = "Hortaliza"
